$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: trim trailing spaces from "Address  " -> "Address"
$ws.Range("J2").Value = "Address"

# Row 3 (Rajesh Kumar) - fill in missing First Name, update mobile/email, fix DOB separator
$ws.Range("L3").Value = "rajesh"
$ws.Range("N3").Value = 9728822877
$ws.Range("S3").Value = "rajeshcdacdemo@cdac.com"

$ws.Range("U3").NumberFormat = "@"
$ws.Range("U3").Value = "10-10-1981"
$ws.Range("U3").NumberFormat = "General"

# Row 4 (Faizal) - fix address typo, fill in missing Last Name, update mobile/email, fix DOB separator
$ws.Range("J4").Value = "ABHOR"
$ws.Range("M4").Value = "aaaa"
$ws.Range("N4").Value = 8338057323
$ws.Range("S4").Value = "faizal1demo@gmail.com"

$ws.Range("U4").NumberFormat = "@"
$ws.Range("U4").Value = "10-10-1981"
$ws.Range("U4").NumberFormat = "General"
